# Updates cryptos list values (Price & Volume(1h) columns) and swaps the
# EthereumClassic / ImmutableX rows, matching the upstream GitHub Actions
# data refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. "27.326.65", "1.006") that
# must stay as text. Mark the cells that are being rewritten as Text format
# first so Excel does not coerce them into floating point numbers (which
# would corrupt values like "1.470" -> 1.47 or introduce FP noise).
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Cell value updates (old -> new), row by row.
$ws.Range("D2").Value = "27.326.65"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").Value = "1.783.59"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "335.55"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "0.3781"
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("D8").Value = "0.3415"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("D9").Value = "48.05"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").Value = "0.07439"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "21.90"
$ws.Range("E13").Value = "  +8.85%  "
$ws.Range("D14").Value = "6.449"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "1.786.69"
$ws.Range("E15").Value = "  +3.47%  "
$ws.Range("D16").Value = "7.007"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "0.00001089"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "0.06637"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "84.31"
$ws.Range("E19").Value = "  +2.99%  "
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "17.28"
$ws.Range("E21").Value = "  +4.51%  "
$ws.Range("D22").Value = "6.438"
$ws.Range("E22").Value = "  +4.67%  "
$ws.Range("D23").Value = "27.307.37"
$ws.Range("E23").Value = "  +3.04%  "
$ws.Range("D24").Value = "12.48"
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("D25").Value = "2.452"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "2.546"
$ws.Range("E26").Value = "  +6.36%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "1.470"
$ws.Range("E27").Value = "  +4.15%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "21.26"
$ws.Range("E28").Value = "  +9.55%  "
$ws.Range("D29").Value = "149.94"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("D30").Value = "1.987.29"
$ws.Range("E30").Value = "  +3.57%  "
$ws.Range("D31").Value = "132.91"
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("D32").Value = "4.061"
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").Value = "6.087"
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("D34").Value = "0.08618"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "13.17"
$ws.Range("E35").Value = "  +3.13%  "
$ws.Range("D36").Value = "1.665"
$ws.Range("E36").Value = "  -1.74%  "
$ws.Range("D37").Value = "0.6857"
$ws.Range("E37").Value = "  +10.51%  "
$ws.Range("D38").Value = "5.407"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "0.06326"
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("D40").Value = "8.776"
$ws.Range("E40").Value = "  +4.36%  "
$ws.Range("D41").Value = "0.02337"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "1.271"
$ws.Range("E43").Value = "  +4.10%  "
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").Value = "1.003"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").Value = "0.6412"
$ws.Range("E46").Value = "  +6.61%  "
$ws.Range("D47").Value = "3.845"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").Value = "2.107"
$ws.Range("E48").Value = "  +3.10%  "
$ws.Range("D49").Value = "129.02"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "78.87"
$ws.Range("E51").Value = "  +2.59%  "
